$d = $word.ActiveDocument

# 1. Remove the "Meta description" paragraph (paragraph 2) entirely, including its
#    paragraph mark, by deleting its full Range.
$metaPara = $d.Paragraphs.Item(2)
$metaPara.Range.Delete()

# 2. Insert a new bold paragraph "Play Action Bank Slot for Free - Read Our Review"
#    right before the final "Please design a feature image..." paragraph, then
#    replace that final paragraph's text with the meta-description text (keeping
#    its existing italic run/formatting untouched).
$found = $d.Content
$found.Find.Execute("Please design a feature image", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$found.Collapse(1)
$insertStart = $found.Start
$found.InsertBefore("Play Action Bank Slot for Free - Read Our Review`r")

$newRunEnd = $insertStart + 49
$newRun = $d.Range($insertStart, $newRunEnd)
$newRun.Font.Bold = $true
$newRun.Font.Italic = $false

$d.Content.Find.Execute("Please design a feature image for the game ""Action Bank"" that fits the following criteria: - The image should be in cartoon style - The image should feature a happy Maya warrior with glasses The Maya warrior in the image should be holding a golden vault with a big smile on his face, representing the potential winnings in the game. The warrior should be wearing a traditional Maya headdress and glasses, emphasizing the modern twist to this classic slot game. The background of the image should be bright and colorful, with bold reel symbols including lucky 7s, Xs, bars, and noughts. This feature image should be eye-catching and capture the fun and excitement of playing Action Bank.", $true, $false, $false, $false, $false, $true, 1, $false, "Looking for a fun, free online slot game? Read our review of Action Bank to see why you should play today.", 2)
